$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Title shape: "COMET DSMC Meeting" + line break  ->  "DSM Report" ---
$titleShape = $s.Shapes.Item("Title 1")
$titleRange = $titleShape.TextFrame.TextRange
$titleRange.Delete()
$titleRange.Text = "DSM Report"

# --- Subtitle shape: drop the grant-number paragraphs, leave it blank ---
$subtitleShape = $s.Shapes.Item("Subtitle 2")
$subtitleRange = $subtitleShape.TextFrame.TextRange
$subtitleRange.Delete()

# --- Remove the logo picture shape entirely ---
$logoShape = $s.Shapes.Item("Picture 4")
$logoShape.Delete()
